# Justify alignment for all non-code body text.
#
# The source document has a handful of body paragraphs (three list items in
# the intro, and two paragraphs that sit right before/after a code-box
# drawing) that were missing the "justify" paragraph alignment that the rest
# of the body text already uses. This sets w:jc="both" (wdAlignParagraphJustify)
# on exactly those paragraphs.

$d = $word.ActiveDocument

# wdAlignParagraphJustify
$wdAlignParagraphJustify = 3

function Set-JustifyByExactText {
    param(
        $Document,
        [string]$Pattern
    )

    $matched = 0
    foreach ($p in $Document.Paragraphs) {
        # Paragraph.Range.Text includes the trailing paragraph mark (CR, and
        # sometimes a cell-mark / BEL); strip those before comparing so the
        # match is exact and not sensitive to the mark character(s).
        $t = $p.Range.Text.TrimEnd([char]13, [char]7)
        if ($t -like $Pattern) {
            $p.Format.Alignment = $wdAlignParagraphJustify
            $matched++
        }
    }
    if ($matched -ne 1) {
        Write-Host "WARNING: pattern matched $matched paragraphs (expected 1): $Pattern"
    }
    return $matched
}

# The three numbered-list ("numId 2") intro bullets.
Set-JustifyByExactText $d "Private and public key generation" | Out-Null
Set-JustifyByExactText $d "Encryption of a text message with the distributed public key " | Out-Null
Set-JustifyByExactText $d "Decoding of the encrypted message using the private key" | Out-Null

# The paragraph holding the "Decryption script" code-box drawing (its own
# visible text is the lead-in sentence before the code block).
Set-JustifyByExactText $d "After receiving the encrypted message (encrypted_message.bin) through Teams, we can decrypt it using the private key. This is done through the following Python code:" | Out-Null

# The closing paragraph after the code box (wildcard avoids fragile
# dependence on the exact smart-quote / accented characters in the quoted
# message).
Set-JustifyByExactText $d "After running the script, ceu.edu could successfully read that the message received indeed translates to*after decryption." | Out-Null
